# Add a new resource entry for "Large Language Model tools for R" as the
# 4th data row (row 4) of the table, pushing the existing rows 4-34 down
# to rows 5-35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 4 (shifts rows 4.. down by one).
$ws.Rows("4:4").Insert()

# Populate the new row with the new resource's data (author, title, link).
$ws.Range("A4").Value = "Luis D. Verde Arregoitia"
$ws.Range("B4").Value = "Large Language Model tools for R"
$ws.Range("C4").Value = "https://luisdva.github.io/llmsr-book/"

# Match the row height used by the other single-line rows in the table.
$ws.Rows("4:4").RowHeight = 17
